# Auto-committed on 2022/05/19 週四 16:27:12.46
#
# Renumber the SEQ column (A12:A21) now that rows 4 and 5 were removed from
# the table (old SEQ 6..15 become 4..13), and switch the "形態" column
# (D9:D17) from center-aligned to left-aligned to match the rest of the
# table's formatting. Also update the sheet's scroll/selection state to
# where the user left off (around row 15, column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# --- Re-number SEQ column (A12:A21) ---------------------------------------
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(14, 1).Value = 6
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(20, 1).Value = 12
$ws.Cells.Item(21, 1).Value = 13

# --- Left-align the "形態" column for rows 9-17 (was centered) ------------
$xlLeft = -4131
$ws.Range("D9:D17").HorizontalAlignment = $xlLeft

# --- Update scroll position / selection to match where the user left off --
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C15").Select()
